# Insert a new data row at row 89 (Ají price entry), shifting all
# subsequent rows (old 89..181) down by one (new 90..182).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(89).Insert()

$ws.Range("A89").Value = 5
$ws.Range("B89").Value = "Macroferia Regional de Talca"
$ws.Range("C89").Value = "Maule"
$ws.Range("D89").Value = 44601
$ws.Range("E89").Value = 7
$ws.Range("F89").Value = 100112021
$ws.Range("G89").Value = "Ají"
$ws.Range("H89").Value = "Americana (o)"
$ws.Range("I89").Value = "Primera"
$ws.Range("J89").Value = 100
$ws.Range("K89").Value = 17000
$ws.Range("L89").Value = 17000
$ws.Range("M89").Value = 17000
$ws.Range("N89").Value = "`$/saco 25 kilos"
$ws.Range("O89").Value = "Región del Maule"
$ws.Range("P89").Value = 680
$ws.Range("Q89").Value = 25
$ws.Range("R89").Value = "Hortaliza"
